$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 308, shifting existing rows 308:403 down to 309:404
$ws.Rows("308:308").Insert()

# Populate the newly inserted row 308 with the new record's values
$ws.Range("A308").Value = 5
$ws.Range("B308").Value = "Macroferia Regional de Talca"
$ws.Range("C308").Value = "Maule"
$ws.Range("D308").Value = 44809
$ws.Range("E308").Value = 7
$ws.Range("F308").Value = 100112023
$ws.Range("G308").Value = "Brócoli"
$ws.Range("H308").Value = "Sin especificar"
$ws.Range("I308").Value = "Primera"
$ws.Range("J308").Value = 2000
$ws.Range("K308").Value = 1200
$ws.Range("L308").Value = 1200
$ws.Range("M308").Value = 1200
$ws.Range("N308").Value = "$/unidad"
$ws.Range("O308").Value = "Región del Maule"
$ws.Range("P308").Value = 1200
$ws.Range("Q308").Value = 1
$ws.Range("R308").Value = "Hortaliza"
